$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 502; this shifts rows 502-537 down to 503-538
# and preserves all their existing cell values/styles.
$ws.Rows.Item(502).Insert()

# Populate the newly inserted row 502 with a fresh data record (new reported
# price for this product). Columns A,B,C,E-L,N-T mirror what was previously
# in row 502 (now shifted to row 503), while D (Fecha) and M (Volumen) carry
# the new values.
$ws.Cells.Item(502, 1).Value = 10
$ws.Cells.Item(502, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(502, 3).Value = "La Araucanía"
$ws.Cells.Item(502, 4).Value = 45013
$ws.Cells.Item(502, 5).Value = 9
$ws.Cells.Item(502, 6).Value = "Fruta"
$ws.Cells.Item(502, 7).Value = 100108
$ws.Cells.Item(502, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(502, 9).Value = 100108002
$ws.Cells.Item(502, 10).Value = "Mango"
$ws.Cells.Item(502, 11).Value = "Sin especificar"
$ws.Cells.Item(502, 12).Value = "Primera"
$ws.Cells.Item(502, 13).Value = 125
$ws.Cells.Item(502, 14).Value = 8000
$ws.Cells.Item(502, 15).Value = 8000
$ws.Cells.Item(502, 16).Value = 8000
$ws.Cells.Item(502, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(502, 18).Value = "Perú"
$ws.Cells.Item(502, 19).Value = 2000
$ws.Cells.Item(502, 20).Value = 4
